$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.95
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 2.63
$ws.Range("K2").Value = 2.2
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("Q2").Value = 1.85
$ws.Range("R2").Value = 1.95
$ws.Range("U2").Value = 1.73
$ws.Range("V2").Value = 2
$ws.Range("X2").Value = 10
$ws.Range("Y2").Value = 9
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 11
$ws.Range("AE2").Value = 15
$ws.Range("AH2").Value = 11
$ws.Range("AI2").Value = 19
$ws.Range("AJ2").Value = 12
$ws.Range("AO2").Value = 11
$ws.Range("AP2").Value = 21
$ws.Range("AS2").Value = 151
$ws.Range("AW2").Value = 5.5
$ws.Range("AX2").Value = 19

# Row 3
$ws.Range("G3").Value = 1.67
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 4.75
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("AA3").Value = 13
$ws.Range("AL3").Value = 41
$ws.Range("AN3").Value = 3.6

# Row 4
$ws.Range("N4").Value = 8.5
$ws.Range("O4").Value = 1.21
$ws.Range("P4").Value = 3.9
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 2.2
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3.05

# Row 5
$ws.Range("G5").Value = 3.4
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 2.2
$ws.Range("J5").Value = 3.75
$ws.Range("L5").Value = 2.88
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.75
$ws.Range("U5").Value = 1.75
$ws.Range("V5").Value = 2
$ws.Range("X5").Value = 17
$ws.Range("Y5").Value = 12
$ws.Range("Z5").Value = 34
$ws.Range("AA5").Value = 26
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 6
$ws.Range("AE5").Value = 13
$ws.Range("AF5").Value = 41
$ws.Range("AG5").Value = 201
$ws.Range("AH5").Value = 8
$ws.Range("AI5").Value = 11
$ws.Range("AJ5").Value = 9
$ws.Range("AK5").Value = 21
$ws.Range("AL5").Value = 17
$ws.Range("AM5").Value = 26
$ws.Range("AN5").Value = 5
$ws.Range("AO5").Value = 19
$ws.Range("AP5").Value = 26
$ws.Range("AQ5").Value = 51
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 201
$ws.Range("AT5").Value = 2.75
$ws.Range("AW5").Value = 4.33
$ws.Range("AX5").Value = 12
$ws.Range("AZ5").Value = 41

# Row 8
$ws.Range("G8").Value = 3.9
$ws.Range("H8").Value = 3.7
$ws.Range("I8").Value = 1.83
$ws.Range("J8").Value = 4.33
$ws.Range("K8").Value = 2.25
$ws.Range("Q8").Value = 1.88
$ws.Range("R8").Value = 1.93
$ws.Range("S8").Value = 1.36
$ws.Range("T8").Value = 3
$ws.Range("W8").Value = 12
$ws.Range("X8").Value = 21
$ws.Range("AP8").Value = 26
$ws.Range("AR8").Value = 81
$ws.Range("AS8").Value = 151
$ws.Range("AT8").Value = 3
$ws.Range("AY8").Value = 19
$ws.Range("BB8").Value = 126

# Row 10
$ws.Range("G10").Value = 2.05
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.2
$ws.Range("J10").Value = 2.75
$ws.Range("K10").Value = 2.05
$ws.Range("L10").Value = 3.75
$ws.Range("M10").Value = 1.03
$ws.Range("N10").Value = 9.5
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3.25
$ws.Range("Q10").Value = 2.05
$ws.Range("R10").Value = 1.75
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("W10").Value = 7.5
$ws.Range("X10").Value = 10
$ws.Range("Y10").Value = 9.5
$ws.Range("Z10").Value = 19
$ws.Range("AA10").Value = 19
$ws.Range("AD10").Value = 6.5
$ws.Range("AE10").Value = 15
$ws.Range("AH10").Value = 9.5
$ws.Range("AI10").Value = 17
$ws.Range("AJ10").Value = 12
$ws.Range("AK10").Value = 34
$ws.Range("AL10").Value = 29
$ws.Range("AM10").Value = 34
$ws.Range("AO10").Value = 12
$ws.Range("AP10").Value = 23
$ws.Range("AQ10").Value = 41
$ws.Range("AT10").Value = 2.63
$ws.Range("AU10").Value = 8
$ws.Range("AW10").Value = 5.5
$ws.Range("AX10").Value = 19
$ws.Range("AZ10").Value = 51
$ws.Range("BA10").Value = 81
$ws.Range("BB10").Value = 300

# Row 11
$ws.Range("Q11").Value = 1.95
$ws.Range("R11").Value = 1.85
$ws.Range("BB11").Value = 400

# Row 13
$ws.Range("I13").Value = 2.15
$ws.Range("J13").Value = 3.75
$ws.Range("O13").Value = 1.36
$ws.Range("P13").Value = 3
$ws.Range("Q13").Value = 2.1
$ws.Range("R13").Value = 1.7
$ws.Range("U13").Value = 1.95
$ws.Range("V13").Value = 1.8
$ws.Range("AG13").Value = 351

Write-Output "Applied 147 cell updates across 8 rows"